$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "SXT"
$ws.Range("B11").Value = 2268
$ws.Range("C11").Value = 2900.696421663087
$ws.Range("D11").Value = 0.6243574749091155
